$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.555.56"
$ws.Range("E2").Value = "  +1.69%  "
$ws.Range("D3").Value = "1.912.77"
$ws.Range("E3").Value = "  +5.36%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.48%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5229"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.73%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3962"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09711"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.153"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.95"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.538"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.74%  "
$ws.Range("E13").Value = "  +3.01%  "
$ws.Range("D14").Value = "1.910.18"
$ws.Range("E14").Value = "  +5.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.558"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.93%  "
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001138"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06650"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.83%  "
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.335"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.43%  "
$ws.Range("D23").Value = "28.643.44"
$ws.Range("E23").Value = "  +1.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.304"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.74%  "
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.698"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +12.08%  "
$ws.Range("B27").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C27").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D27").Value = "2.128.78"
$ws.Range("E27").Value = "  +5.15%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.77%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "158.88"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "129.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.75%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.110"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.36%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1085"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.36%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.753"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.27%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.639"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.55%  "
$ws.Range("B35").Value = "FraxShare"
$ws.Range("C35").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.927"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +11.53%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06771"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.20%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02435"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.00%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.268"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.78%  "
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2232"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.18%  "
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.20%  "
$ws.Range("B41").Value = "InternetComputer(DFINITY)"
$ws.Range("C41").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.096"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6472"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.45%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.191"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.84%  "
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.58"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.28%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6100"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.38%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.756"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.74%  "
$ws.Range("B48").Value = "WEMIXTOKEN"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.286"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.32%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.032"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.74%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "125.39"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.81%  "
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.209"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.54%  "
